$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header M1
$ws.Range("M1").Value = "PAID"

# New donation rows (rows 15-24), columns A-L
$data = @(
    @("REC-1741969487529-984", 14, "2025-03-14T16:24:47.546Z", "Deepak", "Adhikari", "988", "sjahbfkjadhkjfa@gmail.com", "3477712375", "11 alpine ln", "Hicksville", "NY", "11801"),
    @("REC-1741969570079-329", 15, "2025-03-14T16:26:10.085Z", "Deepak", "Adhikari", "988", "sjahbfkjadhkjfa@gmail.com", "3477712375", "11 alpine ln", "Hicksville", "NY", "11801"),
    @("REC-1741969650869-422", 16, "2025-03-14T16:27:30.874Z", "Deepak", "Adhikari", "988", "sjahbfkjadhkjfa@gmail.com", "3477712375", "11 alpine ln", "Hicksville", "NY", "11801"),
    @("REC-1741969739843-247", 17, "2025-03-14T16:28:59.849Z", "Deepak", "Adhikari", "333333", "dadhikari856@gmail.com", "3477712375", "11 alpine ln", "Hicksville", "NY", "11801"),
    @("REC-1741970623658-586", 18, "2025-03-14T16:43:43.665Z", "Deepak", "Adhikari", "2222333", "dadhikari856@gmail.com", "3477712375", "11 alpine ln", "Hicksville", "NY", "11801"),
    @("REC-1741970747227-173", 19, "2025-03-14T16:45:47.232Z", "Deepak", "Adhikari", "2222333", "dadhikari856@gmail.com", "3477712375", "11 alpine ln", "Hicksville", "NY", "11801"),
    @("REC-1741971092861-971", 20, "2025-03-14T16:51:32.867Z", "Deepak", "Adhikari", "2222333", "dadhikari856@gmail.com", "3477712375", "11 alpine ln", "Hicksville", "NY", "11801"),
    @("REC-1741971120418-780", 21, "2025-03-14T16:52:00.420Z", "Deepak", "Adhikari", "6666666666", "dadhikari856@gmail.com", "3477712375", "11 alpine ln", "Hicksville", "NY", "11801"),
    @("REC-1741971739730-822", 22, "2025-03-14T17:02:19.735Z", "Deepak", "Adhikari", "132", "deepak-adhikari@hotmail.com", "8567768105", "11 alpine ln", "Hicksville", "Alabama", "11801"),
    @("REC-1741971785782-752", 23, "2025-03-14T17:03:05.790Z", "Deepak", "Adhikari", "132", "deepak-adhikari@taptap.com", "8567768105", "11 alpine ln", "Hicksville", "Alabama", "11801")
)

$startRow = 15
$endRow = $startRow + $data.Count - 1

# Columns A,C,D,E,F,G,H,I,J,K,L hold text values (numeric-looking strings must
# stay text, matching the existing rows); only column B (Row Number) is numeric.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"
$ws.Range("C$startRow`:L$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $colNum = $c + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowValues[$c]
    }
}
